# SustainabilitySignals.pptx — apply the tracked edit described by the
# commit diff:
#   1) The cached "datetimeFigureOut" field text ("2/15/2026") is bumped to
#      "2/18/2026" everywhere it is cached — the Slide Master's Date
#      placeholder and every Custom Layout's Date placeholder.
#   2) Slide 1's "Rounded Rectangle 6" badge ("Student Project") is removed.
#   3) Slide 11's dual-degree caption textbox drops the trailing
#      "(Woxsen)" qualifier and shrinks to fit the new shorter text.
#   4) Slide 4's "Completeness" label textbox is resized/repositioned.

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached date field text (Master + every Layout) -------
$oldDate = "2/15/2026"
$newDate = "2/18/2026"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2) Slide 1: delete the "Student Project" badge shape ----------------
$slide1 = $p.Slides.Item(1)
for ($i = $slide1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "Rounded Rectangle 6" -and $sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Student Project") {
        $sh.Delete()
    }
}

# --- 3) Slide 11: trim "(Woxsen)" from the dual-degree caption -----------
$slide11 = $p.Slides.Item(11)
for ($i = 1; $i -le $slide11.Shapes.Count; $i++) {
    $sh = $slide11.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "*Dual Degree*") {
        $sh.TextFrame.TextRange.Text = "🎓  Dual Degree: Master in Finance (AMS) + MBA Financial Services "
        $sh.Height = 276999 / 12700
    }
}

# --- 4) Slide 4: reposition/resize the "Completeness" label --------------
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $sh = $slide4.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Completeness") {
        $sh.Left = 1810512 / 12700
        $sh.Width = 1527048 / 12700
        $sh.Height = 353943 / 12700
        $sh.TextFrame.TextRange.Text = "Completeness"
    }
}
